# Generate Report for Handoff
# Updates the localization-status report: the file has moved from
# "In Translation" to "Ready for handoff", and the handoff datetimes
# for both the zh-cn and de-de targets are refreshed.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ---
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Latest Handoff Datetime updates ---
# zh-cn handoff datetime: 2016-03-23 02:34:25 -> 2016-03-23 02:35:41
$zhcn.Range("E2").Value = "2016-03-23 02:35:41"

# de-de handoff datetime: 2016-03-23 02:34:29 -> 2016-03-23 02:35:44
$dede.Range("E2").Value = "2016-03-23 02:35:44"

# Overview "Latest Handoff Date" mirrors the de-de handoff timestamp
$overview.Range("D2").Value = "2016-03-23 02:35:44"
